$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5; existing rows 5-7 shift down to 6-8.
$ws.Rows("5:5").Insert()

# Populate the new row 5 with this week's data.
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = 44846
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 100114007
$ws.Cells.Item(5, 7).Value = "Jengibre"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 18000
$ws.Cells.Item(5, 12).Value = 18000
$ws.Cells.Item(5, 13).Value = 18000
$ws.Cells.Item(5, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(5, 15).Value = "Perú"
$ws.Cells.Item(5, 16).Value = 1385
$ws.Cells.Item(5, 17).Value = 13
$ws.Cells.Item(5, 18).Value = "Hortaliza"
